$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price column stores numeric-looking values as literal TEXT in the
# source data (inline/shared strings, not numeric cells). Force each
# touched Price cell to Text format *before* writing so Excel does not
# auto-convert strings like "235.96" into a Number cell.
$priceCells = @("D2", "D3", "D4", "D5", "D6", "D7", "D8", "D9", "D10", "D11", "D12", "D13", "D14", "D15", "D16", "D17", "D18", "D19", "D20", "D21", "D22", "D23", "D24", "D25", "D26", "D29", "D30", "D31", "D32", "D33", "D34", "D35", "D36", "D37", "D38", "D39", "D40", "D41", "D42", "D43", "D44", "D45", "D46", "D47", "D48", "D49", "D50")
foreach ($addr in $priceCells) {
  $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "30.237.57"
$ws.Range("E2").Value = "  -0.32%  "
$ws.Range("D3").Value = "1.863.17"
$ws.Range("E3").Value = "  -0.39%  "
$ws.Range("D4").Value = "0.9993"
$ws.Range("E4").Value = "  -0.18%  "
$ws.Range("D5").Value = "235.96"
$ws.Range("E5").Value = "  +0.18%  "
$ws.Range("D6").Value = "0.9995"
$ws.Range("E6").Value = "  -0.15%  "
$ws.Range("D7").Value = "0.4702"
$ws.Range("E7").Value = "  +0.31%  "
$ws.Range("D8").Value = "0.2904"
$ws.Range("E8").Value = "  +2.09%  "
$ws.Range("D9").Value = "0.06575"
$ws.Range("E9").Value = "  +0.46%  "
$ws.Range("D10").Value = "21.88"
$ws.Range("E10").Value = "  +1.79%  "
$ws.Range("D11").Value = "0.08007"
$ws.Range("E11").Value = "  +1.45%  "
$ws.Range("D12").Value = "97.80"
$ws.Range("E12").Value = "  -0.30%  "
$ws.Range("D13").Value = "1.860.39"
$ws.Range("E13").Value = "  -0.55%  "
$ws.Range("D14").Value = "5.121"
$ws.Range("E14").Value = "  +0.12%  "
$ws.Range("D15").Value = "0.6794"
$ws.Range("E15").Value = "  +0.48%  "
$ws.Range("D16").Value = "269.62"
$ws.Range("E16").Value = "  -3.08%  "
$ws.Range("D17").Value = "30.218.62"
$ws.Range("E17").Value = "  -0.36%  "
$ws.Range("D18").Value = "13.65"
$ws.Range("E18").Value = "  +7.08%  "
$ws.Range("D19").Value = "0.000007655"
$ws.Range("E19").Value = "  +4.58%  "
$ws.Range("D20").Value = "0.9997"
$ws.Range("E20").Value = "  -0.09%  "
$ws.Range("D21").Value = "2.104.29"
$ws.Range("E21").Value = "  -0.62%  "
$ws.Range("D22").Value = "0.9993"
$ws.Range("E22").Value = "  -0.23%  "
$ws.Range("D23").Value = "5.234"
$ws.Range("E23").Value = "  -4.39%  "
$ws.Range("D24").Value = "6.188"
$ws.Range("E24").Value = "  +0.57%  "
$ws.Range("D25").Value = "167.27"
$ws.Range("E25").Value = "  +1.16%  "
$ws.Range("D26").Value = "9.200"
$ws.Range("E26").Value = "  +0.30%  "
$ws.Range("E27").Value = "  -1.34%  "
$ws.Range("E28").Value = "  +1.00%  "
$ws.Range("D29").Value = "1.372"
$ws.Range("E29").Value = "  -0.82%  "
$ws.Range("D30").Value = "0.09941"
$ws.Range("E30").Value = "  +2.43%  "
$ws.Range("D31").Value = "4.338"
$ws.Range("E31").Value = "  -1.46%  "
$ws.Range("D32").Value = "1.467"
$ws.Range("E32").Value = "  -0.53%  "
$ws.Range("D33").Value = "4.045"
$ws.Range("E33").Value = "  -1.46%  "
$ws.Range("D34").Value = "0.04713"
$ws.Range("E34").Value = "  -0.10%  "
$ws.Range("D35").Value = "1.123"
$ws.Range("E35").Value = "  -0.68%  "
$ws.Range("D36").Value = "0.7037"
$ws.Range("E36").Value = "  -0.43%  "
$ws.Range("D37").Value = "2.717"
$ws.Range("E37").Value = "  -0.42%  "
$ws.Range("D38").Value = "0.01880"
$ws.Range("E38").Value = "  +0.83%  "
$ws.Range("D39").Value = "2.607"
$ws.Range("E39").Value = "  +2.64%  "
$ws.Range("D40").Value = "6.345"
$ws.Range("E40").Value = "  +0.16%  "
$ws.Range("D41").Value = "73.48"
$ws.Range("E41").Value = "  -1.62%  "
$ws.Range("D42").Value = "1.942"
$ws.Range("E42").Value = "  -0.28%  "
$ws.Range("D43").Value = "0.8407"
$ws.Range("E43").Value = "  -1.18%  "
$ws.Range("D44").Value = "103.88"
$ws.Range("E44").Value = "  +0.16%  "
$ws.Range("D45").Value = "0.9990"
$ws.Range("E45").Value = "  -0.18%  "
$ws.Range("D46").Value = "0.4151"
$ws.Range("E46").Value = "  -0.97%  "
$ws.Range("B47").Value = "Aptos"
$ws.Range("C47").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D47").Value = "7.061"
$ws.Range("E47").Value = "  -2.25%  "
$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").Value = "9.162"
$ws.Range("E48").Value = "  -1.06%  "
$ws.Range("D49").Value = "928.80"
$ws.Range("E49").Value = "  -1.10%  "
$ws.Range("D50").Value = "34.16"
$ws.Range("E50").Value = "  -0.34%  "
$ws.Range("E51").Value = "  +0.39%  "
